# Revert unit test coverage: re-add row 39 to each of the 4 worksheets,
# mirroring the structure/format of the existing data rows (e.g. row 38).

$wb = $excel.ActiveWorkbook

# Row data for each sheet, in workbook sheet order:
# DE_LFT_#1, DE_LFT_#2, DE_PLT_#1, DE_PLT_#2
$rowsData = @(
    @{
        A = [double]"45825.43813657408"
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x6C"
        E = "0x14"
        F = 380
        G = [double]"7.598631275147109e+23"
        H = 364
        I = 14
    },
    @{
        A = [double]"45825.43813657408"
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x70"
        E = "0xe"
        F = 380
        G = [double]"5.68432987514711e+23"
        H = 368
        I = 14
    },
    @{
        A = [double]"45825.43813657408"
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x80"
        E = "0x7"
        F = 130
        G = [double]"5.68631262647114e+23"
        H = 128
        I = 7
    },
    @{
        A = [double]"45825.43813657408"
        B = "0x00,0x82"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x80"
        E = "0x3"
        F = 130
        G = [double]"9.85046333984776e+23"
        H = 128
        I = 3
    }
)

for ($i = 1; $i -le 4; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $data = $rowsData[$i - 1]

    $ws.Range("A39").Value = $data.A
    $ws.Range("A39").NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Range("B39").Value = $data.B
    $ws.Range("C39").Value = $data.C
    $ws.Range("D39").Value = $data.D
    $ws.Range("E39").Value = $data.E

    $ws.Range("F39").Value = $data.F
    $ws.Range("G39").Value = $data.G
    $ws.Range("H39").Value = $data.H
    $ws.Range("I39").Value = $data.I
}
